$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("team_codes")

# Rename team names that fbref now reports using longer / different forms.
# (Formulas on the "fixtures" sheet look these names up via INDEX/MATCH,
# so any fixture still referencing an old short name will recalc to #N/A.)
$ws.Range("A19").Value = "Tottenham"
$ws.Range("A11").Value = "Leeds United"
$ws.Range("A16").Value = "Newcastle Utd"
$ws.Range("A17").Value = "Nott'ham Forest"
$ws.Range("A12").Value = "Leicester City"
$ws.Range("A15").Value = "Manchester Utd"
$ws.Range("A14").Value = "Manchester City"

# Leave the cursor where the author left it when they committed this change.
$ws.Range("B4").Select() | Out-Null
